$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values (row -> B,C,D,E,F) as described by the commit "fixed traj read real robot data"
$data = @{
    1 = @(0.9489610112883385, 0.1758320240729323, -0.7249361535492151, 0.6700277771922065, 1.570796383038867)
    2 = @(0.9943554056026752, 0.1769150497322079, -0.721918637663874, 0.6719623020239368, 1.570796377815036)
    3 = @(1.19775238786194, 0.1817677227138239, -0.708398165934528, 0.6806302558279879, 1.570796354408809)
    4 = @(1.483912630473311, 0.1885949728642848, -0.6893761459059412, 0.6928252438556356, 1.570796321478471)
    5 = @(1.687309612732577, 0.1934476458459009, -0.6758556741765951, 0.7014931976596867, 1.570796298072244)
    6 = @(1.732704007046913, 0.1945306715051764, -0.672838158291254, 0.7034277224914169, 1.570796292848413)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
}
